$wb = $excel.ActiveWorkbook

# Rename the "SettlementItems" sheet back to "SettlementInventory"
$ws = $wb.Worksheets.Item("SettlementItems")
$ws.Name = "SettlementInventory"

# Make the renamed sheet the active sheet/tab (moves tabSelected + activeTab)
$ws.Activate()
